$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per latest scrape
$ws.Range('D2').Value = '22.558.81'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '1.578.26'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').Value = '''288.95'
$ws.Range('E6').Value = '  -0.85%  '
$ws.Range('D7').Value = '''0.3698'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').Value = '''48.68'
$ws.Range('E8').Value = '  -2.44%  '
$ws.Range('D9').Value = '''0.3350'
$ws.Range('E9').Value = '  -1.13%  '
$ws.Range('D10').Value = '''1.145'
$ws.Range('E10').Value = '  +0.24%  '
$ws.Range('D11').Value = '''0.07487'
$ws.Range('E11').Value = '  -0.82%  '
$ws.Range('D13').Value = '''21.04'
$ws.Range('E13').Value = '  -1.21%  '
$ws.Range('D14').Value = '''6.011'
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').Value = '''6.971'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').Value = '1.582.18'
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('D17').Value = '''0.00001121'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').Value = '''88.85'
$ws.Range('E18').Value = '  -2.04%  '
$ws.Range('D19').Value = '''0.06768'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').Value = '''6.423'
$ws.Range('E20').Value = '  +1.84%  '
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').Value = '''16.60'
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('D23').Value = '''12.18'
$ws.Range('E23').Value = '  -0.41%  '
$ws.Range('D24').Value = '22.557.90'
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('D26').Value = '''2.600'
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('D27').Value = '''152.73'
$ws.Range('E27').Value = '  +2.37%  '
$ws.Range('D28').Value = '''19.71'
$ws.Range('E28').Value = '  -1.78%  '
$ws.Range('D29').Value = '''5.014'
$ws.Range('E29').Value = '  -0.81%  '
$ws.Range('D30').Value = '''124.51'
$ws.Range('E30').Value = '  -0.61%  '
$ws.Range('D31').Value = '1.755.33'
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('D32').Value = '''1.068'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('D33').Value = '''6.193'
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('D35').Value = '''9.682'
$ws.Range('E35').Value = '  -0.92%  '
$ws.Range('D36').Value = '''0.08339'
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('D37').Value = '''0.02464'
$ws.Range('E37').Value = '  -1.06%  '
$ws.Range('D38').Value = '''0.2271'
$ws.Range('E38').Value = '  -1.44%  '
$ws.Range('D39').Value = '''5.453'
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '''0.06397'
$ws.Range('E40').Value = '  -2.28%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '''1.300'
$ws.Range('E41').Value = '  -4.94%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '''0.6366'
$ws.Range('E42').Value = '  +2.12%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '''11.40'
$ws.Range('E43').Value = '  +0.32%  '
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').Value = '''13.99'
$ws.Range('E45').Value = '  -0.32%  '
$ws.Range('D46').Value = '''0.6197'
$ws.Range('E46').Value = '  +5.72%  '
$ws.Range('D47').Value = '''3.774'
$ws.Range('E47').Value = '  -0.90%  '
$ws.Range('D48').Value = '''2.066'
$ws.Range('E48').Value = '  -0.41%  '
$ws.Range('D49').Value = '''124.99'
$ws.Range('E49').Value = '  -3.26%  '
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('E51').Value = '  -0.72%  '

Write-Output "Applied crypto list update"
